$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 6: B6 and C6 change value (A6 stays the same)
$ws.Range("B6").Value = "Grad der Umsetzung"
$ws.Range("C6").Value = "xxxLevel of compliance"

# Rows 7 through 24 each take on the values that used to live one row below them
$ws.Range("A7").Value = "K_CALCMETH"
$ws.Range("B7").Value = "Berechnungsmethode"
$ws.Range("C7").Value = "Calculation method"

$ws.Range("A8").Value = "K_CRIM"
$ws.Range("B8").Value = "Straftat"
$ws.Range("C8").Value = "Criminal offence"

$ws.Range("A9").Value = "K_CRIMOFF"
$ws.Range("B9").Value = "Straftaten"
$ws.Range("C9").Value = "Criminal offences"

$ws.Range("A10").Value = "K_KREIS"
$ws.Range("B10").Value = "Kreis"
$ws.Range("C10").Value = "County"

$ws.Range("A11").Value = "K_LAENDER"
$ws.Range("B11").Value = "Bundesland"
$ws.Range("C11").Value = "Federal state"

$ws.Range("A12").Value = "K_PM"
$ws.Range("B12").Value = "Feinstaub"
$ws.Range("C12").Value = "Fine particulate matter"

$ws.Range("A13").Value = "K_PRAEV"
$ws.Range("B13").Value = "Art der Maßnahme"
$ws.Range("C13").Value = "XXXArt der Prävention"

$ws.Range("A14").Value = "K_QUALI"
$ws.Range("B14").Value = "Qualifizierung"
$ws.Range("C14").Value = "Qualification"

$ws.Range("A15").Value = "K_RISIKO"
$ws.Range("B15").Value = "Risikolage"
$ws.Range("C15").Value = "XXXRisikolage"

$ws.Range("A16").Value = "K_SEA"
$ws.Range("B16").Value = "Meer"
$ws.Range("C16").Value = "Sea"

$ws.Range("A17").Value = "K_SEKTOREN"
$ws.Range("B17").Value = "Sektoren"
$ws.Range("C17").Value = "XXXSektoren"

$ws.Range("A18").Value = "K_SERIES"
$ws.Range("B18").Value = "Zeitreihe"
$ws.Range("C18").Value = "Time series"

$ws.Range("A19").Value = "K_SEX"
$ws.Range("B19").Value = "Geschlecht"
$ws.Range("C19").Value = "Sex"

$ws.Range("A20").Value = "K_SUBINDEX"
$ws.Range("B20").Value = "Teilindizes"
$ws.Range("C20").Value = "Sub index"

$ws.Range("A21").Value = "K_TARIF"
$ws.Range("B21").Value = "Art des Tarifvertrags"
$ws.Range("C21").Value = "XXXArt der Tarifverträge"

$ws.Range("A22").Value = "K_TYPEAREA"
$ws.Range("B22").Value = "Art der Fläche"
$ws.Range("C22").Value = "Type of area"

$ws.Range("A23").Value = "K_URBAN"
$ws.Range("B23").Value = "Verstädterungsgrad"
$ws.Range("C23").Value = "Degree of urbanisation"

$ws.Range("A24").Value = "K_ZUORDN"
$ws.Range("B24").Value = "Zuordnung"
$ws.Range("C24").Value = "XXXZuordnung"

# Remove the now-duplicate last row (row 25)
$ws.Rows.Item(25).Delete()
